# Add a new leaderboard entry (row 7) to the "Leaderboard" sheet,
# mirroring the existing rows' data/format:
#   A: Participant, B: Course, C: Battery, D: Date, E: Laptime

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Leaderboard")

$ws.Range("A7").Value = "Toby"
$ws.Range("B7").Value = "Trophy Course"
$ws.Range("C7").Value = "2S"
$ws.Range("D7").Value = 46024.5
$ws.Range("E7").Value = 18.809999999999999

# Update the active selection as left by the edit
[void]$ws.Range("D18").Select()
